$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update TF8 (row 9) values: Length 12->6, Width 6->4 ---
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 4

# --- Update TF9 (row 10) values: Length 12->8 (Width stays 6) ---
$ws.Range("B10").Value = 8

# --- Remove the TF10/TF11/TF12 rows (old rows 11-13) ---
$ws.Rows("11:13").Delete()

# --- Relabel the remaining last two symbols to TF11 / TF12 ---
$ws.Range("A9").Value = "TF11"
$ws.Range("A10").Value = "TF12"

# --- Append three new blank (but styled) rows at 17-19 for the new objective prototypes ---
$ws.Range("A10:E10").Copy()
$ws.Range("A17:E19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A17:E19").ClearContents()
$excel.CutCopyMode = $false

# --- Update the active selection to reflect the newly-added rows ---
$ws.Activate()
$ws.Range("A15:E21").Select()
